$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Append new changelog entry for version 1.5.1 (row 20)
$ws.Range("A20").Value = 44368
$ws.Range("B20").Value = "1.5.1"
$ws.Range("C20").Value = "Improvements:`n- show circ supply, total cupply, market cap and corresponding rank also for missing DFI token. NaN-entry handled as Zero and user gets a hint of not correct values."

# Copy formatting (date format, text format, wrap-text format) from the row above
$ws.Range("A19:C19").Copy()
$ws.Range("A20:C20").PasteSpecial(-4122)

# Row height for the new wrapped-text entry
$ws.Rows.Item(20).RowHeight = 45

# Update the selection to reflect where the user ended up after editing
$ws.Range("C21").Select()
